$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Extend the used range with a new column U, copying the formatting that
#    already exists one column to the left (T) for every row so the new
#    column's style matches its row's existing look (borders / fills / etc).
# ---------------------------------------------------------------------------
foreach ($r in 1..6) {
    $ws.Range("T$r").Copy() | Out-Null
    $ws.Range("U$r").PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Row 1 & Row 2: the company-info text that used to live in column M now
#    lives in column N; column M becomes a blank cell (same look as the
#    other blank header cells, e.g. column O).
# ---------------------------------------------------------------------------
$companyName = $ws.Range("M1").Value2
$companyInfo = $ws.Range("M2").Value2

$ws.Range("O1").Copy() | Out-Null
$ws.Range("N1").PasteSpecial(-4122) | Out-Null
$ws.Range("O2").Copy() | Out-Null
$ws.Range("N2").PasteSpecial(-4122) | Out-Null

$ws.Range("O1").Copy() | Out-Null
$ws.Range("M1").PasteSpecial(-4122) | Out-Null
$ws.Range("O2").Copy() | Out-Null
$ws.Range("M2").PasteSpecial(-4122) | Out-Null

$ws.Range("N1").Value = $companyName
$ws.Range("N2").Value = $companyInfo
$ws.Range("M1").Value = $null
$ws.Range("M2").Value = $null

# ---------------------------------------------------------------------------
# 3) Row 4 header: insert a new "igst_amount" header before "cgst_amount",
#    shifting the remaining headers one column to the right (P..T -> Q..U).
#    All of these header cells already share the same style, so only the
#    text needs to move.
# ---------------------------------------------------------------------------
$ws.Range("U4").Value = $ws.Range("T4").Value2
$ws.Range("T4").Value = $ws.Range("S4").Value2
$ws.Range("S4").Value = $ws.Range("R4").Value2
$ws.Range("R4").Value = $ws.Range("Q4").Value2
$ws.Range("Q4").Value = $ws.Range("P4").Value2
$ws.Range("P4").Value = "igst_amount"

# ---------------------------------------------------------------------------
# 4) Row 5 data: same shift as row 4 (new igst_amount value of 2500 before
#    the old cgst_amount value), plus the email/name updates and the new
#    highlight colour on the invoice number cell.
# ---------------------------------------------------------------------------
$ws.Range("U5").Value = $ws.Range("T5").Value2
$ws.Range("T5").Value = $ws.Range("S5").Value2
$ws.Range("S5").Value = $ws.Range("R5").Value2
$ws.Range("R5").Value = $ws.Range("Q5").Value2
$ws.Range("Q5").Value = $ws.Range("P5").Value2
$ws.Range("P5").Value = 2500

$ws.Range("G5").Value = "sachinsakh108@gmail.com"
$ws.Range("I5").Value = "Sachin Prabhu"

$ws.Range("K5").Interior.Color = 65535
